# Auto-generated edit script: updates the cryptos price table
# to match the Thu Sep 12 18:55:00 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.240.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.350.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.53%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.23%  "

# Row 7
$ws.Range("E7").Value = "  +0.58%  "

# Row 8
$ws.Range("E8").Value = "  +4.86%  "

# Row 9
$ws.Range("E9").Value = "  +0.33%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.05%  "

# Row 11
$ws.Range("E11").Value = "  -0.54%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.35%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.771.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.34%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.183.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.85%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.351.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.37%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.83%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "332.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.83%  "

# Row 24
$ws.Range("E24").Value = "  +1.23%  "

# Row 25
$ws.Range("E25").Value = "  +0.53%  "

# Row 26
$ws.Range("E26").Value = "  -3.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.66%  "

# Row 28
$ws.Range("E28").Value = "  +1.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.59%  "

# Row 33
$ws.Range("E33").Value = "  +12.36%  "

# Row 34
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.05%  "

# Row 36
$ws.Range("E36").Value = "  +0.97%  "

# Row 37
$ws.Range("E37").Value = "  -0.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "142.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.54%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.378"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.41%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "288.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.71%  "

# Row 44
$ws.Range("E44").Value = "  +1.54%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.46%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0502"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "

# Row 47
$ws.Range("E47").Value = "  +1.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0219"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "

# Row 49
$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.385"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.31%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "

